# Consumption_Forecast_Historical.xlsx update
#
# The sheet holds a rolling 2-day window of 15-minute (96 per day) consumption
# forecast data (rows 2-193, 96 rows per day). This edit drops the oldest day
# (07.11.2025, previously rows 2-97) and appends a freshly fetched day
# (09.11.2025) after the existing newest day (08.11.2025), which shifts up to
# become the new first day (rows 2-97). Column C (Quarter 1-96) and column D
# (Lookup index formula results) are structurally unchanged by this edit -
# only the Timestamp (A) / Forecasted Consumption (B) values move, and the
# Lookup text strings in the shared-string table rotate accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$quartersPerDay = 96
$firstDataRow = 2          # row holding quarter 1 of the (old) first day
$secondDayFirstRow = 98    # row holding quarter 1 of the (old) second day

# ---------------------------------------------------------------------
# Step 1: shift the old "second day" (rows 98-193) up to become the new
# "first day" (rows 2-97) - this is the old 08.11.2025 data moving to the
# front of the rolling window.
# ---------------------------------------------------------------------
$shiftAB = New-Object 'object[,]' $quartersPerDay,2
$shiftD = New-Object 'object[,]' $quartersPerDay,1

for ($i = 0; $i -lt $quartersPerDay; $i++) {
    $srcRow = $secondDayFirstRow + $i
    $shiftAB[$i,0] = $ws.Cells.Item($srcRow, 1).Value2   # Timestamp (A)
    $shiftAB[$i,1] = $ws.Cells.Item($srcRow, 2).Value2   # Forecasted Consumption (B)
    $shiftD[$i,0]  = $ws.Cells.Item($srcRow, 4).Value2   # Lookup text (D)
}

$ws.Range($ws.Cells.Item($firstDataRow,1), $ws.Cells.Item($firstDataRow + $quartersPerDay - 1,2)).Value = $shiftAB
$ws.Range($ws.Cells.Item($firstDataRow,4), $ws.Cells.Item($firstDataRow + $quartersPerDay - 1,4)).Value = $shiftD

# ---------------------------------------------------------------------
# Step 2: write the freshly retrained forecast for the new day
# (09.11.2025) into rows 98-193, the new "second day" slot.
# ---------------------------------------------------------------------
$newDayConsumption = @(
    5200,5150,5110,5070,5030,5010,4990,4970,4960,4950,
    4950,4950,4950,4950,4950,4960,4990,5010,5030,5060,
    5100,5120,5140,5180,5230,5260,5290,5330,5360,5380,
    5390,5390,5380,5370,5340,5310,5260,5220,5180,5140,
    5090,5060,5040,5020,5000,5000,5000,5010,5040,5060,
    5080,5100,5140,5190,5250,5320,5400,5480,5560,5660,
    5770,5890,6020,6140,6290,6410,6510,6590,6660,6700,
    6710,6710,6700,6690,6670,6640,6590,6530,6490,6420,
    6330,6240,6170,6070,5940,5810,5670,5550,5460,5350,
    5230,5120,5160,5100,5050,5000
)

$newDaySerial = [double]45970   # 09.11.2025 as an Excel date serial
$newDateText = "09.11.2025"

$newAB = New-Object 'object[,]' $quartersPerDay,2
$newD = New-Object 'object[,]' $quartersPerDay,1

for ($i = 0; $i -lt $quartersPerDay; $i++) {
    # Round to 11 decimal places to match the exact binary64 values the source
    # data uses for its quarter-hour timestamps (15-minute fraction-of-day
    # steps), so the new day's serials are bit-identical to how the rest of
    # the sheet was generated.
    $newAB[$i,0] = [Math]::Round($newDaySerial + ($i / [double]$quartersPerDay), 11)
    $newAB[$i,1] = $newDayConsumption[$i]
    $newD[$i,0]  = $newDateText + ($i + 1)
}

$ws.Range($ws.Cells.Item($secondDayFirstRow,1), $ws.Cells.Item($secondDayFirstRow + $quartersPerDay - 1,2)).Value = $newAB
$ws.Range($ws.Cells.Item($secondDayFirstRow,4), $ws.Cells.Item($secondDayFirstRow + $quartersPerDay - 1,4)).Value = $newD
